# (hel-888) description pour indicateur Allocation Ress EJ et ETSAN Excel de comparaison
#
# 1) Scroll the sheet-tab strip so sheet 2 ("Comparaison") is the first
#    displayed tab (workbookView/@firstSheet), while keeping it the active tab.
# 2) Add a new sentence to the "HAPI" indicator description (Lisez-moi!C9):
#    "Pour la campagne et le type d'établissement sélectionnés, cette
#    enveloppe fait partie des 3 premières enveloppes attribuées."
#    right after "HAPI" and before "L'application HAPI « Autres champs »...",
#    preserving the existing bold/plain rich-text runs.

$wb = $excel.ActiveWorkbook

# --- 1) view state -----------------------------------------------------
$compSheet = $wb.Worksheets.Item("Comparaison")
$compSheet.Activate()
$win = $excel.ActiveWindow
$win.FirstVisibleSheet = $compSheet.Index
$win.firstSheet = $compSheet.Index

# --- 2) rich text edit on Lisez-moi!C9 ----------------------------------
$ws = $wb.Worksheets.Item("Lisez-moi")
$cell = $ws.Range("C9")

$bold1 = "Source : "
$hapiWord = "HAPI"
$newSentence = "Pour la campagne et le type d’établissement sélectionnés, cette enveloppe fait partie des 3 premières enveloppes attribuées."
$plain1b = "L’application HAPI « Autres champs » outille le processus de gestion et d’attribution des ressources liées aux enveloppes MIGAC, DAF, USLD, Forfait et FMESPP d’une part, et, le pilotage et l’ordonnancement des dépenses du fonds d’intervention régional (FIR) d’autre part."
$bold2 = "Mode de calcul :"
$plain2 = "Montant des crédits alloués par enveloppes, sous enveloppes et mode délégation, par campagne budgétaire."
$bold3 = "Fréquence :"
$plain3 = " Quotidienne"
$bold4 = "Source(s) : "
$plain4 = "HAPI (HArmonisation et Partage d’Information) - Autres champs."
$plain5 = "Hélios collecte ces données depuis le SI mutualisé des ARS DIAMANT « Décisionnel Inter-ARS pour la Maîtrise et l’Anticipation. », outil décisionnel de pilotage centré sur la régulation de l’offre de soins, abordée sous les aspects des moyens humains, financiers, et productivité."

$plainBlock1 = $hapiWord + "`n" + $newSentence + "`n" + $plain1b + "`n`n"
$plainBlock2 = $plain2 + "`n`n"
$plainBlock3 = $plain3 + "`n`n"
$plainBlock4 = $plain4 + "`n`n" + $plain5

$fullText = $bold1 + $plainBlock1 + $bold2 + "`n" + $plainBlock2 + $bold3 + $plainBlock3 + $bold4 + $plainBlock4

$cell.Value2 = $fullText

# Re-apply per-run formatting (whole-cell assignment resets every run to the
# cell's base font) so the saved file keeps the same rich-text runs as the
# original: four bold "headers" and plain Calibri/12/black text in between.
function Set-RunFont($range, [bool]$bold) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 12
    $range.Font.Color = 0
    $range.Font.Bold = $bold
}

$pos = 1
Set-RunFont $cell.Characters($pos, $bold1.Length) $true
$pos = $pos + $bold1.Length
Set-RunFont $cell.Characters($pos, $plainBlock1.Length) $false
$pos = $pos + $plainBlock1.Length

$bold2nl = $bold2 + "`n"
Set-RunFont $cell.Characters($pos, $bold2nl.Length) $true
$pos = $pos + $bold2nl.Length
Set-RunFont $cell.Characters($pos, $plainBlock2.Length) $false
$pos = $pos + $plainBlock2.Length

Set-RunFont $cell.Characters($pos, $bold3.Length) $true
$pos = $pos + $bold3.Length
Set-RunFont $cell.Characters($pos, $plainBlock3.Length) $false
$pos = $pos + $plainBlock3.Length

Set-RunFont $cell.Characters($pos, $bold4.Length) $true
$pos = $pos + $bold4.Length
Set-RunFont $cell.Characters($pos, $plainBlock4.Length) $false

"done"
